# "oprava exel tabulek done" - fix unit scale (divide by 10) for the
# "Sklon mV" (col B) and "Sklon mBar" (col F) measurement columns on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 2.2955999999999999
$ws.Range("B6").Value = 2.2961
$ws.Range("B7").Value = 2.2963
$ws.Range("B8").Value = 2.2780999999999998
$ws.Range("B9").Value = 2.2963
$ws.Range("B10").Value = 2.2808999999999999

$ws.Range("F5").Value = 0.36795
$ws.Range("F6").Value = 0.36677999999999999
$ws.Range("F7").Value = 0.37136999999999998
$ws.Range("F8").Value = 0.3639
$ws.Range("F9").Value = 0.37136999999999998
$ws.Range("F10").Value = 0.36874000000000001

$wb.Save()

# Update selected cell to reflect where the user last clicked after editing.
$ws.Range("F11").Select()
